$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7 corresponds to year 2025 - update computed values per latest data refresh
$ws.Range("B7").Value = 3664770.1
$ws.Range("C7").Value = -17.51731567573698
$ws.Range("D7").Value = 3224
$ws.Range("E7").Value = 3224
$ws.Range("F7").Value = 1136.715291563276
$ws.Range("G7").Value = 21.16563057062955
